$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A4 value (more precise timestamp)
$ws.Range("A4").Value = 45866.16688976852

# Add new row 5 of data
$ws.Range("A5").Value = 45866.25027959104
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat
$ws.Range("B5").Value = 2025
$ws.Range("C5").Value = 31
$ws.Range("D5").Value = 13.67
$ws.Range("E5").Value = 91.16
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 2.24
$ws.Range("H5").Value = "ENE"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "06:00:24"
